$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 3 is being re-purposed from an "AP" ammo_338_federal entry to a "DMG" entry
# (commit: "Deer hunter as 338 Federal"), with a new highlighted font (theme accent2 color).

# A3: keep same ammo name, but apply the new font style
$a3 = $ws.Range("A3")
$a3.Font.ThemeColor = 6

# B3: change tag from AP to DMG, apply the new font style
$b3 = $ws.Range("B3")
$b3.Value = "DMG"
$b3.Font.ThemeColor = 6

# C3: update base price
$ws.Range("C3").Value = 4000

# E3: apply new font style + number format (value recalculates from formula)
$e3 = $ws.Range("E3")
$e3.Font.ThemeColor = 6
$e3.NumberFormat = "0.000"

# H3: update IRL Joules value
$ws.Range("H3").Value = 10.7

# Match the author's final cursor position in the saved workbook
$ws.Range("E11").Select() | Out-Null
